$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D ("Price") and E ("Volume(1h)") hold numeric- and percentage-
# looking values that are stored as literal text in the workbook (inline
# strings, e.g. "8.760" or "-1.03%"). Plain `.Value = "..."` assignment
# would let Excel auto-convert these to a real number/percentage (losing
# trailing zeros / changing the underlying stored value), so each target
# cell is switched to the Text number format ("@") immediately before its
# new value is written, which keeps it a literal string - matching the
# other (non-numeric) text columns that are unaffected by this.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "326.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.03%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.15%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.709"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "5.92%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08032"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.86%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.030"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.27%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.498"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.70%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.631"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.29%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.14%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9228"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.29%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1258"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.53%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1963"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.50%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.760"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "20.87%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09174"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.61%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03571"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.22%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.1051"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.71%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001289"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.03%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006347"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.92%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.364"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.07%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.24%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.57%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2701"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.49%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04413"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.11%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.15%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004615"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.11%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001188"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-1.00%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02493"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.01%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05328"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.84%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007440"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.26%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009914"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "7.72%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1408"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.47%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.65%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01176"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "9.16%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006677"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.43%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.03%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003037"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-9.19%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002280"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-4.96%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.03%"
